$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 7 new rows at the top of the data (pushes existing rows 2-21 down to 9-28)
$ws.Rows("2:8").Insert()
$ws.Rows("2:8").ClearFormats()

$ws.Cells.Item(2, 1).Value = 2.846464770180838
$ws.Cells.Item(2, 2).Value = -7.078682354518345
$ws.Cells.Item(2, 3).Value = -3.06583663395473
$ws.Cells.Item(3, 1).Value = 2.835772105625697
$ws.Cells.Item(3, 2).Value = -7.043793456895011
$ws.Cells.Item(3, 3).Value = -2.961106973034996
$ws.Cells.Item(4, 1).Value = 2.915691324642726
$ws.Cells.Item(4, 2).Value = -7.111980744770594
$ws.Cells.Item(4, 3).Value = -2.798882663249969
$ws.Cells.Item(5, 1).Value = 3.055751519543784
$ws.Cells.Item(5, 2).Value = -7.220968450818743
$ws.Cells.Item(5, 3).Value = -2.894419597727912
$ws.Cells.Item(6, 1).Value = 2.955562557492938
$ws.Cells.Item(6, 2).Value = -7.134888444628034
$ws.Cells.Item(6, 3).Value = -2.944399050303869
$ws.Cells.Item(7, 1).Value = 2.833344757556915
$ws.Cells.Item(7, 2).Value = -7.338198423385621
$ws.Cells.Item(7, 3).Value = -2.509933024644852
$ws.Cells.Item(8, 1).Value = 3.108331612178258
$ws.Cells.Item(8, 2).Value = -7.027578847748893
$ws.Cells.Item(8, 3).Value = -2.627250722476415

# Append 3 new rows (29-31) after the existing data (which now ends at row 28)
$ws.Cells.Item(29, 1).Value = 2.293728096144541
$ws.Cells.Item(29, 2).Value = -7.460245260170529
$ws.Cells.Item(29, 3).Value = 0.3245020040443973
$ws.Cells.Item(30, 1).Value = 1.675990547452653
$ws.Cells.Item(30, 2).Value = -6.526311159133912
$ws.Cells.Item(30, 3).Value = -0.6771522419793252
$ws.Cells.Item(31, 1).Value = 3.116939672401965
$ws.Cells.Item(31, 2).Value = -5.273013770580299
$ws.Cells.Item(31, 3).Value = -3.711685695818474

